{"js": "const tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\nconst table = tables.items[0];\n\nconst replacements = [\n  { row: 0, col: 0, oldText: \"58\u00d728=1624\", newText: \"17\u00d765=1105\" },\n  { row: 0, col: 1, oldText: \"51\u00d781=4131\", newText: \"28\u00d785=2380\" },\n  { row: 0, col: 2, oldText: \"94\u00d748=4512\", newText: \"84\u00d755=4620\" },\n  { row: 0, col: 3, oldText: \"69\u00d772=4968\", newText: \"39\u00d774=2886\" },\n  { row: 0, col: 4, oldText: \"14\u00d761=854\", newText: \"47\u00d756=2632\" },\n  { row: 4, col: 0, oldText: \"71\u00d794=6674\", newText: \"37\u00d762=2294\" },\n  { row: 4, col: 1, oldText: \"47\u00d743=2021\", newText: \"15\u00d792=1380\" },\n  { row: 4, col: 2, oldText: \"68\u00d714=952\", newText: \"22\u00d792=2024\" },\n  { row: 4, col: 3, oldText: \"75\u00d774=5550\", newText: \"16\u00d758=928\" },\n  { row: 4, col: 4, oldText: \"42\u00d785=3570\", newText: \"31\u00d784=2604\" },\n  { row: 9, col: 0, oldText: \"71\u00d751=3621\", newText: \"72\u00d721=1512\" },\n  { row: 9, col: 1, oldText: \"71\u00d784=5964\", newText: \"16\u00d745=720\" },\n  { row: 9, col: 2, oldText: \"18\u00d715=270\", newText: \"29\u00d797=2813\" },\n  { row: 9, col: 3, oldText: \"56\u00d718=1008\", newText: \"41\u00d737=1517\" },\n  { row: 9, col: 4, oldText: \"98\u00d724=2352\", newText: \"52\u00d718=936\" },\n  { row: 14, col: 0, oldText: \"85\u00d769=5865\", newText: \"66\u00d758=3828\" },\n  { row: 14, col: 1, oldText: \"76\u00d713=988\", newText: \"39\u00d747=1833\" },\n  { row: 14, col: 2, oldText: \"59\u00d744=2596\", newText: \"32\u00d795=3040\" },\n  { row: 14, col: 3, oldText: \"61\u00d732=1952\", newText: \"20\u00d769=1380\" },\n  { row: 14, col: 4, oldText: \"83\u00d763=5229\", newText: \"33\u00d718=594\" },\n  { row: 19, col: 0, oldText: \"92\u00d726=2392\", newText: \"54\u00d756=3024\" },\n  { row: 19, col: 1, oldText: \"42\u00d752=2184\", newText: \"25\u00d730=750\" },\n  { row: 19, col: 2, oldText: \"95\u00d722=2090\", newText: \"34\u00d747=1598\" },\n  { row: 19, col: 3, oldText: \"14\u00d761=854\", newText: \"27\u00d726=702\" },\n  { row: 19, col: 4, oldText: \"62\u00d771=4402\", newText: \"88\u00d792=8096\" },\n];\n\n// Get all target cells first.\nconst cells = replacements.map(r => table.getCell(r.row, r.col));\n\n// For each cell, search its body for the old text and queue the replace.\nconst searchResults = [];\nfor (let i = 0; i < replacements.length; i++) {\n  const results = cells[i].body.search(replacements[i].oldText, { matchCase: true, matchWildcards: false });\n  results.load(\"items\");\n  searchResults.push(results);\n}\nawait context.sync();\n\nfor (let i = 0; i < replacements.length; i++) {\n  const items = searchResults[i].items;\n  if (items.length > 0) {\n    items[0].insertText(replacements[i].newText, \"Replace\");\n  } else {\n    // Fallback: replace the whole cell body text, preserving as much as possible.\n    cells[i].body.insertText(replacements[i].newText, \"Replace\");\n  }\n}\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n$t = $d.Tables.Item(1)\n\n$replacements = @(\n  @{ Row = 1; Col = 1; OldText = \"58\u00d728=1624\"; NewText = \"17\u00d765=1105\" },\n  @{ Row = 1; Col = 2; OldText = \"51\u00d781=4131\"; NewText = \"28\u00d785=2380\" },\n  @{ Row = 1; Col = 3; OldText = \"94\u00d748=4512\"; NewText = \"84\u00d755=4620\" },\n  @{ Row = 1; Col = 4; OldText = \"69\u00d772=4968\"; NewText = \"39\u00d774=2886\" },\n  @{ Row = 1; Col = 5; OldText = \"14\u00d761=854\"; NewText = \"47\u00d756=2632\" },\n  @{ Row = 5; Col = 1; OldText = \"71\u00d794=6674\"; NewText = \"37\u00d762=2294\" },\n  @{ Row = 5; Col = 2; OldText = \"47\u00d743=2021\"; NewText = \"15\u00d792=1380\" },\n  @{ Row = 5; Col = 3; OldText = \"68\u00d714=952\"; NewText = \"22\u00d792=2024\" },\n  @{ Row = 5; Col = 4; OldText = \"75\u00d774=5550\"; NewText = \"16\u00d758=928\" },\n  @{ Row = 5; Col = 5; OldText = \"42\u00d785=3570\"; NewText = \"31\u00d784=2604\" },\n  @{ Row = 10; Col = 1; OldText = \"71\u00d751=3621\"; NewText = \"72\u00d721=1512\" },\n  @{ Row = 10; Col = 2; OldText = \"71\u00d784=5964\"; NewText = \"16\u00d745=720\" },\n  @{ Row = 10; Col = 3; OldText = \"18\u00d715=270\"; NewText = \"29\u00d797=2813\" },\n  @{ Row = 10; Col = 4; OldText = \"56\u00d718=1008\"; NewText = \"41\u00d737=1517\" },\n  @{ Row = 10; Col = 5; OldText = \"98\u00d724=2352\"; NewText = \"52\u00d718=936\" },\n  @{ Row = 15; Col = 1; OldText = \"85\u00d769=5865\"; NewText = \"66\u00d758=3828\" },\n  @{ Row = 15; Col = 2; OldText = \"76\u00d713=988\"; NewText = \"39\u00d747=1833\" },\n  @{ Row = 15; Col = 3; OldText = \"59\u00d744=2596\"; NewText = \"32\u00d795=3040\" },\n  @{ Row = 15; Col = 4; OldText = \"61\u00d732=1952\"; NewText = \"20\u00d769=1380\" },\n  @{ Row = 15; Col = 5; OldText = \"83\u00d763=5229\"; NewText = \"33\u00d718=594\" },\n  @{ Row = 20; Col = 1; OldText = \"92\u00d726=2392\"; NewText = \"54\u00d756=3024\" },\n  @{ Row = 20; Col = 2; OldText = \"42\u00d752=2184\"; NewText = \"25\u00d730=750\" },\n  @{ Row = 20; Col = 3; OldText = \"95\u00d722=2090\"; NewText = \"34\u00d747=1598\" },\n  @{ Row = 20; Col = 4; OldText = \"14\u00d761=854\"; NewText = \"27\u00d726=702\" },\n  @{ Row = 20; Col = 5; OldText = \"62\u00d771=4402\"; NewText = \"88\u00d792=8096\" },\n)\n\nforeach ($rep in $replacements) {\n    $cell = $t.Cell($rep.Row, $rep.Col)\n    $cellRange = $cell.Range\n    # Cell.Range.Text includes a trailing cell-mark (CR+BEL); strip it for a sanity check.\n    $current = $cellRange.Text.TrimEnd([char]13, [char]7)\n    if ($current -ne $rep.OldText) {\n        Write-Output \"Warning: cell ($($rep.Row),$($rep.Col)) expected '$($rep.OldText)' but found '$current'\"\n    }\n    # Assigning .Text replaces only the run content; paragraph/run formatting\n    # (font, size, alignment) on the existing run is preserved.\n    $cellRange.Text = $rep.NewText\n}\n"}
